$d = $word.ActiveDocument

$replacements = @(
    @("86-67=", "57-16="),
    @("52+16=", "41+28="),
    @("28+51=", "49-37="),
    @("70+22=", "60-24="),
    @("19-5=", "24+70="),
    @("88-20=", "83-30="),
    @("13-6=", "16+21="),
    @("62-5=", "87-43="),
    @("36+7=", "15+16="),
    @("86-16=", "66-49="),
    @("51-3=", "16+82="),
    @("74-52=", "4+38="),
    @("41+55=", "68+24="),
    @("54-45=", "34+29="),
    @("12-4=", "52-16="),
    @("45+35=", "33+66="),
    @("8+15=", "96-65="),
    @("72-62=", "4+52="),
    @("78-4=", "98-91="),
    @("93-65=", "41-39="),
    @("14-11=", "13+68="),
    @("27+66=", "71-54="),
    @("28+26=", "89-14="),
    @("61+6=", "18+67="),
    @("20+57=", "99-59="),
    @("18+43=", "43+53="),
    @("40-29=", "13+20="),
    @("23-17=", "79+4="),
    @("64-10=", "34-34="),
    @("78+14=", "37+28="),
    @("67+28=", "20-9="),
    @("32-5=", "48+17="),
    @("2+13=", "33-2="),
    @("86-62=", "6+16="),
    @("86-13=", "23+1="),
    @("76-36=", "85-10="),
    @("93-22=", "8+32="),
    @("17+6=", "67+17="),
    @("62-4=", "9+39="),
    @("9+60=", "91-28="),
    @("3+53=", "4+73="),
    @("62-0=", "52+40="),
    @("58-29=", "9+12="),
    @("85-69=", "11+83="),
    @("38+6=", "28-1="),
    @("85-72=", "20+60="),
    @("91-85=", "25+41="),
    @("65+27=", "61+36="),
    @("30-2=", "0+14="),
    @("96+1=", "72-50="),
    @("95-70=", "89-34="),
    @("33+17=", "59+26="),
    @("31+38=", "9+80="),
    @("13+51=", "10+27="),
    @("92+0=", "29-25="),
    @("93-42=", "87-48="),
    @("87-39=", "67-31="),
    @("50+15=", "11+71="),
    @("69-9=", "32-9="),
    @("24+2=", "78-37="),
    @("46-28=", "52-32="),
    @("96-9=", "15+56="),
    @("95-32=", "9+61="),
    @("70+8=", "56+7="),
    @("68-18=", "40+21="),
    @("91-87=", "26+4="),
    @("75-25=", "31+52="),
    @("44+34=", "98-24="),
    @("61+3=", "80-25="),
    @("42-39=", "92-34="),
    @("67-43=", "39+37="),
    @("24+11=", "36-31="),
    @("80-29=", "29+68="),
    @("23-21=", "31-20="),
    @("51+46=", "7+29="),
    @("39+58=", "13+76="),
    @("52-42=", "72-24="),
    @("38+56=", "19-7="),
    @("85+3=", "38+60="),
    @("43+11=", "99-57="),
    @("30+29=", "88-42="),
    @("77-8=", "8+63="),
    @("65+5=", "44-39="),
    @("85-22=", "45+28="),
    @("65-45=", "12+69="),
    @("65+33=", "88+6="),
    @("22+39=", "23+29="),
    @("1+66=", "68-41="),
    @("1+21=", "17+36="),
    @("20+3=", "82+2="),
    @("66+23=", "90-85="),
    @("21+11=", "56-46="),
    @("97-38=", "62-10="),
    @("74-42=", "66-9="),
    @("89-29=", "36-12="),
    @("1+42=", "30-26="),
    @("59-15=", "90+6="),
    @("26+6=", "38-31="),
    @("51-19=", "15-9="),
    @("90-75=", "74-66="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
